$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'51.820.74"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "'  -0.09%  "
$ws.Cells.Item(2, 5).Style = "Normal"

$ws.Cells.Item(3, 4).Value = "'2.897.35"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "'  +3.24%  "
$ws.Cells.Item(3, 5).Style = "Normal"

$ws.Cells.Item(4, 5).Value = "'  +0.18%  "
$ws.Cells.Item(4, 5).Style = "Normal"

$ws.Cells.Item(5, 4).Value = "'352.23"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "'  +0.33%  "
$ws.Cells.Item(5, 5).Style = "Normal"

$ws.Cells.Item(6, 4).Value = "'112.63"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "'  -0.12%  "
$ws.Cells.Item(6, 5).Style = "Normal"

$ws.Cells.Item(7, 4).Value = "'0.556"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "'  -0.30%  "
$ws.Cells.Item(7, 5).Style = "Normal"

$ws.Cells.Item(8, 5).Value = "'  +0.19%  "
$ws.Cells.Item(8, 5).Style = "Normal"

$ws.Cells.Item(9, 4).Value = "'0.618"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "'  -0.52%  "
$ws.Cells.Item(9, 5).Style = "Normal"

$ws.Cells.Item(10, 4).Value = "'39.29"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "'  -2.66%  "
$ws.Cells.Item(10, 5).Style = "Normal"

$ws.Cells.Item(11, 4).Value = "'0.0872"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "'  +3.89%  "
$ws.Cells.Item(11, 5).Style = "Normal"

$ws.Cells.Item(12, 5).Value = "'  +0.77%  "
$ws.Cells.Item(12, 5).Style = "Normal"

$ws.Cells.Item(13, 4).Value = "'19.89"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "'  +0.05%  "
$ws.Cells.Item(13, 5).Style = "Normal"

$ws.Cells.Item(14, 2).Value = "Polkadot"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(14, 4).Value = "'7.69"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "'  -1.14%  "
$ws.Cells.Item(14, 5).Style = "Normal"

$ws.Cells.Item(15, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(15, 4).Value = "'3.365.66"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "'  +3.77%  "
$ws.Cells.Item(15, 5).Style = "Normal"

$ws.Cells.Item(16, 4).Value = "'2.900.69"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "'  +3.89%  "
$ws.Cells.Item(16, 5).Style = "Normal"

$ws.Cells.Item(17, 4).Value = "'0.978"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "'  +0.86%  "
$ws.Cells.Item(17, 5).Style = "Normal"

$ws.Cells.Item(18, 4).Value = "'51.903.72"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "'  +0.13%  "
$ws.Cells.Item(18, 5).Style = "Normal"

$ws.Cells.Item(19, 2).Value = "Uniswap"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(19, 4).Value = "'7.54"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "'  -1.43%  "
$ws.Cells.Item(19, 5).Style = "Normal"

$ws.Cells.Item(20, 2).Value = "ImmutableX"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(20, 4).Value = "'3.28"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "'  -3.00%  "
$ws.Cells.Item(20, 5).Style = "Normal"

$ws.Cells.Item(21, 4).Value = "'14.05"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "'  +3.47%  "
$ws.Cells.Item(21, 5).Style = "Normal"

$ws.Cells.Item(22, 4).Value = "'0.0₃0974"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "'  -0.24%  "
$ws.Cells.Item(22, 5).Style = "Normal"

$ws.Cells.Item(23, 4).Value = "'70.85"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "'  +0.56%  "
$ws.Cells.Item(23, 5).Style = "Normal"

$ws.Cells.Item(24, 4).Value = "'267.83"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "'  -0.55%  "
$ws.Cells.Item(24, 5).Style = "Normal"

$ws.Cells.Item(25, 4).Value = "'2.79"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "'  +1.28%  "
$ws.Cells.Item(25, 5).Style = "Normal"

$ws.Cells.Item(26, 4).Value = "'0.177"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "'  +9.17%  "
$ws.Cells.Item(26, 5).Style = "Normal"

$ws.Cells.Item(27, 4).Value = "'26.66"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "'  +1.83%  "
$ws.Cells.Item(27, 5).Style = "Normal"

$ws.Cells.Item(28, 4).Value = "'1.00"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "'  -0.06%  "
$ws.Cells.Item(28, 5).Style = "Normal"

$ws.Cells.Item(29, 4).Value = "'6.91"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "'  +12.25%  "
$ws.Cells.Item(29, 5).Style = "Normal"

$ws.Cells.Item(30, 4).Value = "'10.57"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "'  +0.96%  "
$ws.Cells.Item(30, 5).Style = "Normal"

$ws.Cells.Item(31, 4).Value = "'0.102"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "'  +12.18%  "
$ws.Cells.Item(31, 5).Style = "Normal"

$ws.Cells.Item(32, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(32, 4).Value = "'36.86"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "'  -5.04%  "
$ws.Cells.Item(32, 5).Style = "Normal"

$ws.Cells.Item(33, 2).Value = "RenderToken"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(33, 4).Value = "'5.98"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "'  +5.68%  "
$ws.Cells.Item(33, 5).Style = "Normal"

$ws.Cells.Item(34, 2).Value = "OKB"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(34, 4).Value = "'52.91"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "'  +0.30%  "
$ws.Cells.Item(34, 5).Style = "Normal"

$ws.Cells.Item(35, 2).Value = "Toncoin"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(35, 4).Value = "'2.10"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "'  -7.59%  "
$ws.Cells.Item(35, 5).Style = "Normal"

$ws.Cells.Item(36, 4).Value = "'0.0449"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "'  -1.09%  "
$ws.Cells.Item(36, 5).Style = "Normal"

$ws.Cells.Item(37, 4).Value = "'0.999"
$ws.Cells.Item(37, 4).Style = "Normal"

$ws.Cells.Item(38, 4).Value = "'3.33"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "'  +4.94%  "
$ws.Cells.Item(38, 5).Style = "Normal"

$ws.Cells.Item(39, 4).Value = "'18.61"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "'  -1.89%  "
$ws.Cells.Item(39, 5).Style = "Normal"

$ws.Cells.Item(40, 4).Value = "'2.03"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "'  +0.62%  "
$ws.Cells.Item(40, 5).Style = "Normal"

$ws.Cells.Item(41, 4).Value = "'2.68"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "'  +6.15%  "
$ws.Cells.Item(41, 5).Style = "Normal"

$ws.Cells.Item(42, 4).Value = "'0.116"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "'  -0.17%  "
$ws.Cells.Item(42, 5).Style = "Normal"

$ws.Cells.Item(43, 4).Value = "'22.93"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "'  +3.43%  "
$ws.Cells.Item(43, 5).Style = "Normal"

$ws.Cells.Item(44, 4).Value = "'2.18"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "'  -2.37%  "
$ws.Cells.Item(44, 5).Style = "Normal"

$ws.Cells.Item(45, 2).Value = "ApeXProtocol"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Cells.Item(45, 4).Value = "'2.51"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "'  +1.78%  "
$ws.Cells.Item(45, 5).Style = "Normal"

$ws.Cells.Item(46, 2).Value = "Maker"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(46, 4).Value = "'2.177.30"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "'  +2.33%  "
$ws.Cells.Item(46, 5).Style = "Normal"

$ws.Cells.Item(47, 2).Value = "NEARProtocol"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(47, 4).Value = "'3.49"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "'  -1.56%  "
$ws.Cells.Item(47, 5).Style = "Normal"

$ws.Cells.Item(48, 4).Value = "'114.60"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "'  -6.08%  "
$ws.Cells.Item(48, 5).Style = "Normal"

$ws.Cells.Item(49, 5).Value = "'  +12.30%  "
$ws.Cells.Item(49, 5).Style = "Normal"

$ws.Cells.Item(50, 4).Value = "'0.0342"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "'  +6.96%  "
$ws.Cells.Item(50, 5).Style = "Normal"

$ws.Cells.Item(51, 4).Value = "'0.940"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "'  -5.03%  "
$ws.Cells.Item(51, 5).Style = "Normal"
